$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Top summary rows (1-indexed) get their values swapped for the ones that
# used to live in the tab-separated breakdown rows near the bottom, and a
# couple of independent numeric corrections.
$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "3893"
$t.Cell(5, 1).Range.Text  = "0.00002"
$t.Cell(7, 1).Range.Text  = "0.00016"
$t.Cell(8, 1).Range.Text  = "0.00005"
$t.Cell(9, 1).Range.Text  = "0.00027"
$t.Cell(10, 1).Range.Text = "0.00031"
$t.Cell(11, 1).Range.Text = "0.00038"
$t.Cell(12, 1).Range.Text = "0.71379"

# The bottom three rows used to hold a full tab-separated breakdown; they now
# collapse down to just their leading value (which moved up to the top rows).
$t.Cell(44, 1).Range.Text = "99.91"
$t.Cell(45, 1).Range.Text = "0.71"
$t.Cell(46, 1).Range.Text = "775"
